$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 475
$ws1.Range("F3").Value = 5695
$ws1.Range("F5").Value = 71
$ws1.Range("F6").Value = 96
$ws1.Range("F9").Value = 542

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 475
$ws4.Range("F3").Value = 5695
$ws4.Range("F6").Value = 71
$ws4.Range("F7").Value = 96
$ws4.Range("F11").Value = 542
